$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '27.041.95'
$ws.Range("E2").Value = '  -0.39%  '
$ws.Range("D3").Value = '1.620.42'
$ws.Range("E3").Value = '  -1.07%  '
$ws.Range("D5").Value = '213.75'
$ws.Range("E5").Value = '  -1.38%  '
$ws.Range("E6").Value = '  -1.22%  '
$ws.Range("E7").Value = '  +0.04%  '
$ws.Range("E8").Value = '  +0.26%  '
$ws.Range("E10").Value = '  -0.41%  '
$ws.Range("D11").Value = '0.0841'
$ws.Range("E11").Value = '  -0.77%  '
$ws.Range("D12").Value = '1.848.67'
$ws.Range("E12").Value = '  -0.98%  '
$ws.Range("D13").Value = '1.617.59'
$ws.Range("E13").Value = '  -1.16%  '
$ws.Range("E14").Value = '  -0.18%  '
$ws.Range("E15").Value = '  -0.47%  '
$ws.Range("D16").Value = '27.035.34'
$ws.Range("E16").Value = '  -0.40%  '
$ws.Range("D17").Value = '64.39'
$ws.Range("E17").Value = '  -3.22%  '
$ws.Range("E18").Value = '  -0.43%  '
$ws.Range("D19").Value = '214.82'
$ws.Range("E19").Value = '  -0.95%  '
$ws.Range("E20").Value = '  +0.00%  '
$ws.Range("D21").Value = '6.82'
$ws.Range("E21").Value = '  +0.11%  '
$ws.Range("D22").Value = '4.33'
$ws.Range("E22").Value = '  -1.73%  '
$ws.Range("D23").Value = '2.36'
$ws.Range("E23").Value = '  -5.98%  '
$ws.Range("E24").Value = '  -1.15%  '
$ws.Range("D25").Value = '147.39'
$ws.Range("E25").Value = '  +0.34%  '
$ws.Range("D26").Value = '7.43'
$ws.Range("E26").Value = '  +0.40%  '
$ws.Range("E27").Value = '  +0.18%  '
$ws.Range("D28").Value = '0.114'
$ws.Range("E28").Value = '  -3.87%  '
$ws.Range("E29").Value = '  -0.95%  '
$ws.Range("D30").Value = '0.0511'
$ws.Range("E30").Value = '  +0.97%  '
$ws.Range("E31").Value = '  -0.91%  '
$ws.Range("E32").Value = '  -1.87%  '
$ws.Range("D33").Value = '0.725'
$ws.Range("E33").Value = '  +33.26%  '
$ws.Range("E34").Value = '  -0.82%  '
$ws.Range("D35").Value = '1.335.17'
$ws.Range("E35").Value = '  +2.67%  '
$ws.Range("E36").Value = '  -0.66%  '
$ws.Range("E37").Value = '  -0.26%  '
$ws.Range("E38").Value = '  -0.02%  '
$ws.Range("D39").Value = '0.838'
$ws.Range("E39").Value = '  -1.81%  '
$ws.Range("E41").Value = '  -0.52%  '
$ws.Range("D42").Value = '0.795'
$ws.Range("E42").Value = '  -1.60%  '
$ws.Range("D43").Value = '5.37'
$ws.Range("E43").Value = '  +1.40%  '
$ws.Range("D44").Value = '63.87'
$ws.Range("E44").Value = '  +3.61%  '
$ws.Range("D45").Value = '1.759.97'
$ws.Range("E45").Value = '  -0.98%  '
$ws.Range("D46").Value = '90.01'
$ws.Range("E46").Value = '  -1.40%  '
$ws.Range("E47").Value = '  +2.26%  '
$ws.Range("D48").Value = '0.848'
$ws.Range("E48").Value = '  +27.00%  '
$ws.Range("B49").Value = 'Cronos'
$ws.Range("C49").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D49").Value = '0.0511'
$ws.Range("E49").Value = '  -0.24%  '
$ws.Range("B50").Value = 'Algorand'
$ws.Range("C50").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D50").Value = '0.0994'
$ws.Range("E50").Value = '  +3.98%  '
$ws.Range("D51").Value = '7.58'
$ws.Range("E51").Value = '  -1.03%  '
